$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "26.140.43"
$ws.Cells.Item(2, 5).Value2 = "  -4.32%  "

$ws.Cells.Item(3, 4).Value2 = "1.650.90"
$ws.Cells.Item(3, 5).Value2 = "  -3.66%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value2 = "1.009"
$ws.Cells.Item(4, 5).Value2 = "  +0.21%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "216.11"
$ws.Cells.Item(5, 5).Value2 = "  -3.89%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = "0.5124"
$ws.Cells.Item(6, 5).Value2 = "  -3.24%  "

$ws.Cells.Item(7, 5).Value2 = "  +0.21%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = "0.2591"
$ws.Cells.Item(8, 5).Value2 = "  -2.22%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = "0.06434"
$ws.Cells.Item(9, 5).Value2 = "  -3.67%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = "19.76"
$ws.Cells.Item(10, 5).Value2 = "  -5.37%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = "0.07831"
$ws.Cells.Item(11, 5).Value2 = "  +1.65%  "

$ws.Cells.Item(12, 4).Value2 = "1.661.50"

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = "4.287"
$ws.Cells.Item(13, 5).Value2 = "  -4.32%  "

$ws.Cells.Item(14, 4).Value2 = "1.882.72"
$ws.Cells.Item(14, 5).Value2 = "  -3.46%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = "0.5501"
$ws.Cells.Item(15, 5).Value2 = "  -4.91%  "

$ws.Cells.Item(16, 4).Value2 = "0.0₅8002"
$ws.Cells.Item(16, 5).Value2 = "  -2.17%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = "64.04"
$ws.Cells.Item(17, 5).Value2 = "  -5.40%  "

$ws.Cells.Item(18, 4).Value2 = "26.187.57"
$ws.Cells.Item(18, 5).Value2 = "  -4.28%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = "1.008"
$ws.Cells.Item(19, 5).Value2 = "  +0.14%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = "209.47"
$ws.Cells.Item(20, 5).Value2 = "  -4.70%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = "4.400"
$ws.Cells.Item(21, 5).Value2 = "  -5.23%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = "10.07"
$ws.Cells.Item(22, 5).Value2 = "  -3.38%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = "6.008"
$ws.Cells.Item(23, 5).Value2 = "  -0.22%  "

$ws.Cells.Item(24, 5).Value2 = "  +0.21%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = "144.68"
$ws.Cells.Item(25, 5).Value2 = "  -0.55%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "1.788"
$ws.Cells.Item(26, 5).Value2 = "  +4.76%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "0.1172"
$ws.Cells.Item(27, 5).Value2 = "  -2.95%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = "6.982"
$ws.Cells.Item(28, 5).Value2 = "  -3.52%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = "15.83"
$ws.Cells.Item(29, 5).Value2 = "  -2.36%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = "0.05087"
$ws.Cells.Item(30, 5).Value2 = "  -5.40%  "

$ws.Cells.Item(31, 5).Value2 = "  -4.19%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "3.353"
$ws.Cells.Item(32, 5).Value2 = "  -3.58%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = "3.230"
$ws.Cells.Item(33, 5).Value2 = "  -4.74%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "1.554"
$ws.Cells.Item(34, 5).Value2 = "  -5.01%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "2.735"
$ws.Cells.Item(35, 5).Value2 = "  -4.11%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "2.356"
$ws.Cells.Item(36, 5).Value2 = "  -1.76%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = "0.9174"
$ws.Cells.Item(37, 5).Value2 = "  -3.58%  "

$ws.Cells.Item(38, 2).Value2 = "Maker"
$ws.Cells.Item(38, 3).Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(38, 4).Value2 = "1.167.83"
$ws.Cells.Item(38, 5).Value2 = "  +0.81%  "

$ws.Cells.Item(39, 2).Value2 = "ImmutableX"
$ws.Cells.Item(39, 3).Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = "0.5698"
$ws.Cells.Item(39, 5).Value2 = "  -3.27%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "0.01584"
$ws.Cells.Item(40, 5).Value2 = "  -3.92%  "

$ws.Cells.Item(41, 2).Value2 = "PaxDollar"
$ws.Cells.Item(41, 3).Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = "1.008"
$ws.Cells.Item(41, 5).Value2 = "  +0.14%  "

$ws.Cells.Item(42, 2).Value2 = "mCoin"
$ws.Cells.Item(42, 3).Value2 = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "2.561"
$ws.Cells.Item(42, 5).Value2 = "  -0.70%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = "5.678"
$ws.Cells.Item(43, 5).Value2 = "  -2.70%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = "0.8285"
$ws.Cells.Item(44, 5).Value2 = "  -1.43%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = "100.13"
$ws.Cells.Item(45, 5).Value2 = "  -0.81%  "

$ws.Cells.Item(46, 4).Value2 = "1.797.40"
$ws.Cells.Item(46, 5).Value2 = "  -3.25%  "

$ws.Cells.Item(47, 4).Value2 = "0.0₈111"
$ws.Cells.Item(47, 5).Value2 = "  -6.18%  "

$ws.Cells.Item(48, 5).Value2 = "  -0.47%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "55.30"
$ws.Cells.Item(49, 5).Value2 = "  -4.16%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = "1.003"
$ws.Cells.Item(50, 5).Value2 = "  -0.22%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = "7.865"
$ws.Cells.Item(51, 5).Value2 = "  -3.20%  "
